# Apply new message-level prediction results to the "predidx" (D) and
# "pred_name" (E) columns for the rows that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Row = 3;   D = "[1, 0, 0, 1, 0, 0, 0]"; E = "['Normal', 'ParamViolation']" },
    @{ Row = 9;   D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 11;  D = "[1, 0, 1, 0, 1, 0, 0]"; E = "['Normal', 'HardwareFault', 'RegulationViolation']" },
    @{ Row = 16;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 27;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 29;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 35;  D = "[0, 0, 1, 0, 0, 0, 0]"; E = "['HardwareFault']" },
    @{ Row = 54;  D = "[0, 0, 1, 0, 0, 0, 0]"; E = "['HardwareFault']" },
    @{ Row = 61;  D = "[0, 0, 0, 0, 0, 0, 1]"; E = "['SoftwareFault']" },
    @{ Row = 69;  D = "[1, 1, 0, 0, 0, 1, 0]"; E = "['Normal', 'SurroundingEnvironment', 'CommunicationIssue']" },
    @{ Row = 73;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 80;  D = "[1, 0, 1, 0, 0, 0, 0]"; E = "['Normal', 'HardwareFault']" },
    @{ Row = 82;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 83;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 84;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 97;  D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" },
    @{ Row = 116; D = "[1, 0, 0, 0, 0, 0, 0]"; E = "['Normal']" }
)

foreach ($change in $changes) {
    $ws.Cells.Item($change.Row, 4).Value = $change.D
    $ws.Cells.Item($change.Row, 5).Value = $change.E
}
